$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53: No Accounting for Waste / Enchanted Electrum Ink
$ws.Range("H53").Value = 381
$ws.Range("I53").Value = 114.75
$ws.Range("J53").Value = 469.75
$ws.Range("K53").Value = 114.75
$ws.Range("L53").Value = 469.75
$ws.Range("M53").Value = 522.25
$ws.Range("N53").Value = -1743.75

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 1461.2222
$ws.Range("I98").Value = 1493.2858
$ws.Range("J98").Value = 1349
$ws.Range("K98").Value = 1493.2858
$ws.Range("L98").Value = 1349
$ws.Range("M98").Value = 4.714199999999892
$ws.Range("N98").Value = -4345

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 1461.2222
$ws.Range("I122").Value = 1493.2858
$ws.Range("J122").Value = 1349
$ws.Range("K122").Value = 4479.857400000001
$ws.Range("L122").Value = 4047
$ws.Range("M122").Value = -2029.857400000001
$ws.Range("N122").Value = -8947

# Row 130: Technically Still Magic / Ophiotauroskin Magitek Codex
$ws.Range("H130").Value = 58571.43
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 58571.43
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 58571.43
$ws.Range("N130").Value = -68611.42999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 4125.6
$ws.Range("I32").Value = 4311.0977
$ws.Range("J32").Value = 2224.25
$ws.Range("K32").Value = 4311.0977
$ws.Range("L32").Value = 2224.25
$ws.Range("M32").Value = -4024.0977
$ws.Range("N32").Value = -2798.25

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1369.7435
$ws.Range("I74").Value = 1390.6333
$ws.Range("J74").Value = 1300.1111
$ws.Range("K74").Value = 1390.6333
$ws.Range("L74").Value = 1300.1111
$ws.Range("M74").Value = -516.6333
$ws.Range("N74").Value = -3048.1111

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1369.7435
$ws.Range("I77").Value = 1390.6333
$ws.Range("J77").Value = 1300.1111
$ws.Range("K77").Value = 6953.166499999999
$ws.Range("L77").Value = 6500.5555
$ws.Range("M77").Value = -2585.166499999999
$ws.Range("N77").Value = -15236.5555

# Row 111: Hedging Bets / Deepgold Surcoat of Maiming
$ws.Range("H111").Value = 22792.5
$ws.Range("I111").Value = 585
$ws.Range("J111").Value = 45000
$ws.Range("K111").Value = 585
$ws.Range("L111").Value = 45000
$ws.Range("M111").Value = 3505
$ws.Range("N111").Value = -53180

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 5483.1924
$ws.Range("I122").Value = 5368.154
$ws.Range("J122").Value = 5598.231
$ws.Range("K122").Value = 16104.462
$ws.Range("L122").Value = 16794.693
$ws.Range("M122").Value = -13654.462
$ws.Range("N122").Value = -21694.693

# Row 135: Forgiveness for My Shins / Ruthenium Sabatons of Fending
$ws.Range("H135").Value = 69999.55
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 69999.55
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 69999.55
$ws.Range("N135").Value = -80139.55

$ws = $wb.Worksheets.Item("BSM")
# Row 64: With Bearings Straight / Mythrite Nugget
$ws.Range("H64").Value = 1232.1111
$ws.Range("I64").Value = 498.33334
$ws.Range("J64").Value = 1599
$ws.Range("K64").Value = 498.33334
$ws.Range("L64").Value = 1599
$ws.Range("M64").Value = -273.33334
$ws.Range("N64").Value = -2049

# Row 67: Bearing the Brunt (L) / Mythrite Nugget
$ws.Range("H67").Value = 1232.1111
$ws.Range("I67").Value = 498.33334
$ws.Range("J67").Value = 1599
$ws.Range("K67").Value = 498.33334
$ws.Range("L67").Value = 1599
$ws.Range("M67").Value = 281.66666
$ws.Range("N67").Value = -3159

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 12501199
$ws.Range("I86").Value = 15626197
$ws.Range("J86").Value = 1207.25
$ws.Range("K86").Value = 15626197
$ws.Range("L86").Value = 1207.25
$ws.Range("M86").Value = -15625074
$ws.Range("N86").Value = -3453.25

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 12501199
$ws.Range("I89").Value = 15626197
$ws.Range("J89").Value = 1207.25
$ws.Range("K89").Value = 78130985
$ws.Range("L89").Value = 6036.25
$ws.Range("M89").Value = -78125369
$ws.Range("N89").Value = -17268.25

# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 3282.2632
$ws.Range("I94").Value = 2798.5386
$ws.Range("J94").Value = 4330.3335
$ws.Range("K94").Value = 2798.5386
$ws.Range("L94").Value = 4330.3335
$ws.Range("M94").Value = -2347.5386
$ws.Range("N94").Value = -5232.3335

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 10082.9
$ws.Range("I107").Value = 10908.154
$ws.Range("J107").Value = 8550.286
$ws.Range("K107").Value = 10908.154
$ws.Range("L107").Value = 8550.286
$ws.Range("M107").Value = -8988.154
$ws.Range("N107").Value = -12390.286

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 289.7857
$ws.Range("I22").Value = 216.8
$ws.Range("J22").Value = 472.25
$ws.Range("K22").Value = 216.8
$ws.Range("L22").Value = 472.25
$ws.Range("M22").Value = 133.2
$ws.Range("N22").Value = -1172.25

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 1497.1111
$ws.Range("I107").Value = 1395.125
$ws.Range("J107").Value = 2313
$ws.Range("K107").Value = 1395.125
$ws.Range("L107").Value = 2313
$ws.Range("M107").Value = 524.875
$ws.Range("N107").Value = -6153

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 3437.7856
$ws.Range("I132").Value = 2511
$ws.Range("J132").Value = 8998.5
$ws.Range("K132").Value = 7533
$ws.Range("L132").Value = 26995.5
$ws.Range("M132").Value = -5003
$ws.Range("N132").Value = -32055.5

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 9746.5
$ws.Range("I134").Value = 9746.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 29239.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -26704.5

$ws = $wb.Worksheets.Item("GSM")
# Row 21: Forever 21K / Brass Ring
$ws.Range("H21").Value = 264000
$ws.Range("I21").Value = 264000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 264000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -263827

# Row 30: Dog Tags Are for Dogs / Brass Ring
$ws.Range("H30").Value = 264000
$ws.Range("I30").Value = 264000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 264000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -263895

# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 132233.22
$ws.Range("I70").Value = 191686.83
$ws.Range("J70").Value = 13326
$ws.Range("K70").Value = 191686.83
$ws.Range("L70").Value = 13326
$ws.Range("M70").Value = -191416.83
$ws.Range("N70").Value = -13866

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 132233.22
$ws.Range("I73").Value = 191686.83
$ws.Range("J73").Value = 13326
$ws.Range("K73").Value = 191686.83
$ws.Range("L73").Value = 13326
$ws.Range("M73").Value = -190750.83
$ws.Range("N73").Value = -15198

# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 17358.334
$ws.Range("I97").Value = 794.5
$ws.Range("J97").Value = 25640.25
$ws.Range("K97").Value = 794.5
$ws.Range("L97").Value = 25640.25
$ws.Range("M97").Value = -298.5
$ws.Range("N97").Value = -26632.25

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 6157.6924
$ws.Range("I102").Value = 4703.9
$ws.Range("J102").Value = 11003.667
$ws.Range("K102").Value = 4703.9
$ws.Range("L102").Value = 11003.667
$ws.Range("M102").Value = -3081.9
$ws.Range("N102").Value = -14247.667

$ws = $wb.Worksheets.Item("LTW")
# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 49116.84
$ws.Range("I100").Value = 62495.844
$ws.Range("J100").Value = 6750
$ws.Range("K100").Value = 62495.844
$ws.Range("L100").Value = 6750
$ws.Range("M100").Value = -61954.844
$ws.Range("N100").Value = -7832

$ws = $wb.Worksheets.Item("WVR")
# Row 54: No Country for Cold Men / Woolen Tights
$ws.Range("H54").Value = 20437.25
$ws.Range("I54").Value = 25000
$ws.Range("J54").Value = 19785.428
$ws.Range("K54").Value = 25000
$ws.Range("L54").Value = 19785.428
$ws.Range("M54").Value = -24480
$ws.Range("N54").Value = -20825.428

# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 1949.0667
$ws.Range("I81").Value = 1363
$ws.Range("J81").Value = 4293.3335
$ws.Range("K81").Value = 2726
$ws.Range("L81").Value = 8586.666999999999
$ws.Range("M81").Value = -1665
$ws.Range("N81").Value = -10708.667

# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 1949.0667
$ws.Range("I84").Value = 1363
$ws.Range("J84").Value = 4293.3335
$ws.Range("K84").Value = 13630
$ws.Range("L84").Value = 42933.335
$ws.Range("M84").Value = -8326
$ws.Range("N84").Value = -53541.335

# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 867.9
$ws.Range("I107").Value = 867.9
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2603.7
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -683.6999999999998

# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 479.6111
$ws.Range("I113").Value = 456.66666
$ws.Range("J113").Value = 594.3333
$ws.Range("K113").Value = 1369.99998
$ws.Range("L113").Value = 1782.9999
$ws.Range("M113").Value = 800.0000199999999
$ws.Range("N113").Value = -6122.9999

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 3948.4
$ws.Range("I132").Value = 2880.6667
$ws.Range("J132").Value = 5550
$ws.Range("K132").Value = 8642.000100000001
$ws.Range("L132").Value = 16650
$ws.Range("M132").Value = -6112.000100000001
$ws.Range("N132").Value = -21710

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1432.8889
$ws.Range("I136").Value = 1432.8889
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4298.6667
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1748.6667
